# Update "想去人数" (want-to-go count) values that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) rows 2-7
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value = 369
$wsExhibit.Range("F3").Value = 788
$wsExhibit.Range("F4").Value = 272
$wsExhibit.Range("F5").Value = 861
$wsExhibit.Range("F6").Value = 2118
$wsExhibit.Range("F7").Value = 190

# Sheet "全部类型" (4th sheet) rows 2,3,4,7,8,10 (mirrors the "展览" rows plus others)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 369
$wsAll.Range("F3").Value = 788
$wsAll.Range("F4").Value = 272
$wsAll.Range("F7").Value = 861
$wsAll.Range("F8").Value = 2118
$wsAll.Range("F10").Value = 190
